$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'29.345.48"
$ws.Range("E2").Value = "  +1.74%  "

$ws.Range("D3").Value = "'1.845.96"
$ws.Range("E3").Value = "  +0.73%  "

$ws.Range("D4").Value = "'1.0000"
$ws.Range("E4").Value = "  +0.00%  "

$ws.Range("D5").Value = "'244.53"
$ws.Range("E5").Value = "  +0.07%  "

$ws.Range("D6").Value = "'0.6877"
$ws.Range("E6").Value = "  -0.84%  "

$ws.Range("E7").Value = "  +0.07%  "

$ws.Range("D8").Value = "'0.3030"
$ws.Range("E8").Value = "  -0.48%  "

$ws.Range("D9").Value = "'0.07549"
$ws.Range("E9").Value = "  -1.69%  "

$ws.Range("D10").Value = "'23.31"
$ws.Range("E10").Value = "  +0.50%  "

$ws.Range("D11").Value = "'0.07671"
$ws.Range("E11").Value = "  -1.55%  "

$ws.Range("B12").Value = "WrappedEther"
$ws.Range("C12").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D12").Value = "'1.850.36"
$ws.Range("E12").Value = "  +0.95%  "

$ws.Range("B13").Value = "Polkadot"
$ws.Range("C13").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("D13").Value = "'5.094"
$ws.Range("E13").Value = "  +0.16%  "

$ws.Range("D14").Value = "'0.6866"
$ws.Range("E14").Value = "  +1.18%  "

$ws.Range("D15").Value = "'89.44"
$ws.Range("E15").Value = "  -3.87%  "

$ws.Range("D16").Value = "'6.289"
$ws.Range("E16").Value = "  -2.35%  "

$ws.Range("D17").Value = "'29.341.48"
$ws.Range("E17").Value = "  +1.65%  "

$ws.Range("D18").Value = "'0.000008206"
$ws.Range("E18").Value = "  -0.71%  "

$ws.Range("D19").Value = "'2.093.91"
$ws.Range("E19").Value = "  +0.61%  "

$ws.Range("D20").Value = "'234.71"
$ws.Range("E20").Value = "  -3.10%  "

$ws.Range("D21").Value = "'12.60"
$ws.Range("E21").Value = "  -0.57%  "

$ws.Range("E22").Value = "  +0.02%  "

$ws.Range("D23").Value = "'7.566"
$ws.Range("E23").Value = "  +1.73%  "

$ws.Range("D24").Value = "'1.001"
$ws.Range("E24").Value = "  -0.02%  "

$ws.Range("D25").Value = "'0.1463"
$ws.Range("E25").Value = "  -1.36%  "

$ws.Range("B26").Value = "Monero"
$ws.Range("C26").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D26").Value = "'159.80"
$ws.Range("E26").Value = "  +0.43%  "

$ws.Range("B27").Value = "Cosmos"
$ws.Range("C27").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("D27").Value = "'8.867"
$ws.Range("E27").Value = "  +1.09%  "

$ws.Range("D28").Value = "'18.08"
$ws.Range("E28").Value = "  -0.86%  "

$ws.Range("D29").Value = "'1.522"
$ws.Range("E29").Value = "  -1.08%  "

$ws.Range("D30").Value = "'4.232"
$ws.Range("E30").Value = "  +0.44%  "

$ws.Range("D31").Value = "'4.123"
$ws.Range("E31").Value = "  -0.56%  "

$ws.Range("D32").Value = "'1.197"
$ws.Range("E32").Value = "  +0.99%  "

$ws.Range("D33").Value = "'0.05184"
$ws.Range("E33").Value = "  +1.97%  "

$ws.Range("D34").Value = "'0.7686"
$ws.Range("E34").Value = "  -0.54%  "

$ws.Range("D35").Value = "'1.866"
$ws.Range("E35").Value = "  +0.81%  "

$ws.Range("D36").Value = "'1.140"
$ws.Range("E36").Value = "  +0.09%  "

$ws.Range("D37").Value = "'2.676"
$ws.Range("E37").Value = "  -0.63%  "

$ws.Range("D38").Value = "'1.301.71"
$ws.Range("E38").Value = "  +5.66%  "

$ws.Range("D39").Value = "'0.01852"
$ws.Range("E39").Value = "  +0.29%  "

$ws.Range("D40").Value = "'2.705"
$ws.Range("E40").Value = "  +0.31%  "

$ws.Range("D41").Value = "'0.9424"
$ws.Range("E41").Value = "  -1.03%  "

$ws.Range("D42").Value = "'105.47"
$ws.Range("E42").Value = "  -2.07%  "

$ws.Range("B43").Value = "PaxDollar"
$ws.Range("C43").Value = "https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp"
$ws.Range("D43").Value = "'1.000"
$ws.Range("E43").Value = "  +0.03%  "

$ws.Range("B44").Value = "FraxShare"
$ws.Range("C44").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D44").Value = "'5.729"
$ws.Range("E44").Value = "  -3.72%  "

$ws.Range("D45").Value = "'9.717"
$ws.Range("E45").Value = "  +1.18%  "

$ws.Range("D46").Value = "'1.993.16"
$ws.Range("E46").Value = "  +0.83%  "

$ws.Range("B47").Value = "BabyDogeCoin"
$ws.Range("C47").Value = "https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge"
$ws.Range("D47").Value = "'0.00000000124"
$ws.Range("E47").Value = "  +5.60%  "

$ws.Range("B48").Value = "Mantle"
$ws.Range("C48").Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$ws.Range("D48").Value = "'0.5214"
$ws.Range("E48").Value = "  +1.06%  "

$ws.Range("B49").Value = "RenderToken"
$ws.Range("C49").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D49").Value = "'1.769"
$ws.Range("E49").Value = "  +1.81%  "

$ws.Range("D50").Value = "'62.96"
$ws.Range("E50").Value = "  -1.23%  "

$ws.Range("D51").Value = "'0.05931"
